# Apply the "updated Gemini response" edit to the Estimated Bills sheet:
# - Insert two new rate-plan columns (Saver's Choice 12, Solar Buyback Saver 36)
#   between the existing "e-Saver" and "Flex Forward" columns (old I shifts to K).
# - Rename a couple of existing plan headers.
# - Replace the usage/estimate figures for every month + the Total row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column I; this shifts the old
# "Flex Forward" column (I) to K, and the new I:J inherit the header/
# blank-row styling automatically from their neighbours.
$ws.Columns("I:J").Insert()

# [row, column, new value] for every cell that actually changes content.
$data = @(
    @(1, 5, ' Simple Rate 12'),
    @(1, 8, ' e-Saver 12'),
    @(1, 9, ' Saver''s Choice 12'),
    @(1, 10, ' Solar Buyback Saver 36'),
    @(2, 2, 100),
    @(2, 3, 29.85),
    @(2, 4, 37.93),
    @(2, 5, 22.88),
    @(2, 6, 33.08),
    @(2, 7, 29.44),
    @(2, 8, 22.5),
    @(2, 9, 24.7),
    @(2, 10, 27.57),
    @(2, 11, 22.68),
    @(3, 2, 200),
    @(3, 3, 49.55),
    @(3, 4, 65.63),
    @(3, 5, 35.68),
    @(3, 6, 55.98),
    @(3, 7, 48.74),
    @(3, 8, 37.9),
    @(3, 9, 39.3),
    @(3, 10, 40.08),
    @(3, 11, 35.28),
    @(4, 2, 300),
    @(4, 3, 69.25),
    @(4, 4, 93.33),
    @(4, 5, 48.48),
    @(4, 6, 78.88),
    @(4, 7, 68.04000000000001),
    @(4, 8, 53.3),
    @(4, 9, 53.9),
    @(4, 10, 52.58),
    @(4, 11, 47.88),
    @(5, 2, 400),
    @(5, 3, 88.95),
    @(5, 4, 121.03),
    @(5, 5, 61.28),
    @(5, 6, 101.78),
    @(5, 7, 87.34),
    @(5, 8, 68.7),
    @(5, 9, 68.5),
    @(5, 10, 65.08),
    @(5, 11, 60.48),
    @(6, 2, 500),
    @(6, 3, 108.65),
    @(6, 4, 148.73),
    @(6, 5, 74.08),
    @(6, 6, 124.68),
    @(6, 7, 106.64),
    @(6, 8, 84.09999999999999),
    @(6, 9, 83.09999999999999),
    @(6, 10, 77.58),
    @(6, 11, 73.08),
    @(7, 2, 600),
    @(7, 3, 128.35),
    @(7, 4, 176.43),
    @(7, 5, 86.88),
    @(7, 6, 147.58),
    @(7, 7, 125.94),
    @(7, 8, 99.5),
    @(7, 9, 97.7),
    @(7, 10, 90.08),
    @(7, 11, 85.68000000000001),
    @(8, 2, 700),
    @(8, 3, 148.05),
    @(8, 4, 204.13),
    @(8, 5, 99.68000000000001),
    @(8, 6, 170.48),
    @(8, 7, 145.24),
    @(8, 8, 114.9),
    @(8, 9, 112.3),
    @(8, 10, 102.58),
    @(8, 11, 98.28),
    @(9, 2, 800),
    @(9, 3, 167.75),
    @(9, 4, 231.83),
    @(9, 5, 112.48),
    @(9, 6, 193.38),
    @(9, 7, 164.54),
    @(9, 8, 130.3),
    @(9, 9, 126.9),
    @(9, 10, 115.08),
    @(9, 11, 110.88),
    @(10, 2, 900),
    @(10, 3, 187.45),
    @(10, 4, 259.53),
    @(10, 5, 125.28),
    @(10, 6, 216.28),
    @(10, 7, 183.84),
    @(10, 8, 145.7),
    @(10, 9, 141.5),
    @(10, 10, 127.58),
    @(10, 11, 123.48),
    @(11, 5, 138.08),
    @(11, 8, 161.1),
    @(11, 9, 156.1),
    @(11, 10, 140.07),
    @(12, 2, 1100),
    @(12, 3, 226.85),
    @(12, 4, 314.93),
    @(12, 5, 150.88),
    @(12, 6, 262.08),
    @(12, 7, 222.44),
    @(12, 8, 171.1),
    @(12, 9, 170.7),
    @(12, 10, 152.57),
    @(12, 11, 148.68),
    @(13, 2, 1200),
    @(13, 3, 246.55),
    @(13, 4, 342.63),
    @(13, 5, 163.68),
    @(13, 6, 284.98),
    @(13, 7, 241.74),
    @(13, 8, 181.1),
    @(13, 9, 185.3),
    @(13, 10, 165.07),
    @(13, 11, 161.28),
    @(14, 3, 1658.4),
    @(14, 4, 2283.36),
    @(14, 5, 1119.36),
    @(14, 6, 1908.36),
    @(14, 7, 1627.08),
    @(14, 8, 1270.2),
    @(14, 9, 1260),
    @(14, 10, 1155.92),
    @(14, 11, 1103.76)
)

foreach ($item in $data) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}

# The blank trailer row (15) keeps the same bold/empty styling across the
# whole row; make the two newly-inserted cells explicit empty text like
# their neighbours instead of leaving them completely value-less.
$ws.Cells.Item(15, 9).Value = ""
$ws.Cells.Item(15, 10).Value = ""
